# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets: row 3 (F3) 80 -> 81, row 4 (F4) 51 -> 52.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 81
    $ws.Range("F4").Value = 52
}
